$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.4004911058505485
$ws.Range("D2").Value = 0.1775387999900744
$ws.Range("E2").Value = 0.1699669624180657
$ws.Range("F2").Value = 1.597962881869364
$ws.Range("G2").Value = 0.9467706423634326
$ws.Range("H2").Value = 1.002803097962662
$ws.Range("J2").Value = 0.210103688300407
$ws.Range("K2").Value = 1.426289054253061
$ws.Range("L2").Value = 0.1545022747935629
$ws.Range("M2").Value = 0.4732552272329684
$ws.Range("O2").Value = 3.942727425474814

$ws.Range("C3").Value = 0.3984529761517024
$ws.Range("D3").Value = 0.1750985820047859
$ws.Range("E3").Value = 0.1704545377486326
$ws.Range("F3").Value = 1.612544587557991
$ws.Range("G3").Value = 0.9584404617549751
$ws.Range("H3").Value = 1.013325435510744
$ws.Range("J3").Value = 0.2120668672593542
$ws.Range("K3").Value = 1.27599483284439
$ws.Range("L3").Value = 0.1553946711331839
$ws.Range("M3").Value = 0.4436045141356857
$ws.Range("O3").Value = 3.98918062432368

$ws.Range("C4").Value = 0.3973803935531208
$ws.Range("D4").Value = 0.1736494436727369
$ws.Range("E4").Value = 0.1708149865162714
$ws.Range("F4").Value = 1.622409667606682
$ws.Range("G4").Value = 0.9662759061860626
$ws.Range("H4").Value = 1.020266409887462
$ws.Range("J4").Value = 0.2133549660341938
$ws.Range("K4").Value = 1.183477214139288
$ws.Range("L4").Value = 0.1559813545475546
$ws.Range("M4").Value = 0.4254287797841769
$ws.Range("O4").Value = 4.020118264284591

$ws.Range("C5").Value = 0.3969883787473378
$ws.Range("D5").Value = 0.1730713518486482
$ws.Range("E5").Value = 0.1709772600650581
$ws.Range("F5").Value = 1.626659036820406
$ws.Range("G5").Value = 0.9696372449163846
$ws.Range("H5").Value = 1.023215687712295
$ws.Range("J5").Value = 0.2139006863152861
$ws.Range("K5").Value = 1.145718826017458
$ws.Range("L5").Value = 0.1562301941095008
$ws.Range("M5").Value = 0.4180301142712821
$ws.Range("O5").Value = 4.033332580731738

$ws.Range("C6").Value = 0.3969260101654015
$ws.Range("D6").Value = 0.1729761141000878
$ws.Range("E6").Value = 0.1710051357235152
$ws.Range("F6").Value = 1.627378487915308
$ws.Range("G6").Value = 0.9702055542387527
$ws.Range("H6").Value = 1.023712708117216
$ws.Range("J6").Value = 0.2139925600188288
$ws.Range("K6").Value = 1.13944573016488
$ws.Range("L6").Value = 0.1562721038087993
$ws.Range("M6").Value = 0.4168020771813588
$ws.Range("O6").Value = 4.035563457786836

$ws.Range("C7").Value = 0.397374924075848
$ws.Range("D7").Value = 0.1736415968318639
$ws.Range("E7").Value = 0.1708171126498677
$ws.Range("F7").Value = 1.622466047818165
$ws.Range("G7").Value = 0.9663205570452078
$ws.Range("H7").Value = 1.020305695843405
$ws.Range("J7").Value = 0.213362241532657
$ws.Range("K7").Value = 1.182968217286316
$ws.Range("L7").Value = 0.1559846709379862
$ws.Range("M7").Value = 0.4253289652642778
$ws.Range("O7").Value = 4.020294020009857

$ws.Range("C8").Value = 0.3997513001753248
$ws.Range("D8").Value = 0.1766872558888224
$ws.Range("E8").Value = 0.1701224198465852
$ws.Range("F8").Value = 1.602801384188837
$ws.Range("G8").Value = 0.9506552047688572
$ws.Range("H8").Value = 1.006331543646738
$ws.Range("J8").Value = 0.2107634369187199
$ws.Range("K8").Value = 1.374518246349567
$ws.Range("L8").Value = 0.1548019414662321
$ws.Range("M8").Value = 0.4630257816551975
$ws.Range("O8").Value = 3.958243013201553

$ws.Range("C9").Value = 0.405826583711459
$ws.Range("D9").Value = 0.1830467611996909
$ws.Range("E9").Value = 0.1692435109495776
$ws.Range("F9").Value = 1.571475646132384
$ws.Range("G9").Value = 0.9252606736441038
$ws.Range("H9").Value = 0.9827375699472043
$ws.Range("J9").Value = 0.2063227007782373
$ws.Range("K9").Value = 1.748170772065237
$ws.Range("L9").Value = 0.1527892691894106
$ws.Range("M9").Value = 0.5371653855701908
$ws.Range("O9").Value = 3.855739807802195

$ws.Range("C10").Value = 0.411148763645258
$ws.Range("D10").Value = 0.1879514165152045
$ws.Range("E10").Value = 0.1688909357086317
$ws.Range("F10").Value = 1.552874819108133
$ws.Range("G10").Value = 0.9098607704214174
$ws.Range("H10").Value = 0.967723792050208
$ws.Range("J10").Value = 0.203458751806977
$ws.Range("K10").Value = 2.021380018939794
$ws.Range("L10").Value = 0.1514964100918341
$ws.Range("M10").Value = 0.5917448029399424
$ws.Range("O10").Value = 3.792143255187

$ws.Range("C11").Value = 0.4137556125314461
$ws.Range("D11").Value = 0.1902324020370685
$ws.Range("E11").Value = 0.1687938838198662
$ws.Range("F11").Value = 1.545371923682254
$ws.Range("G11").Value = 0.9035646676945959
$ws.Range("H11").Value = 0.9613972403257378
$ws.Range("J11").Value = 0.2022422296213087
$ws.Range("K11").Value = 2.145363871828351
$ws.Range("L11").Value = 0.1509483791369224
$ws.Range("M11").Value = 0.6165936478054874
$ws.Range("O11").Value = 3.765759297288099

$ws.Range("C12").Value = 0.4147693842236038
$ws.Range("D12").Value = 0.1911032455550838
$ws.Range("E12").Value = 0.1687662127892189
$ws.Range("F12").Value = 1.542668684095752
$ws.Range("G12").Value = 0.9012827070004903
$ws.Range("H12").Value = 0.9590738995148911
$ws.Range("J12").Value = 0.201793960809848
$ws.Range("K12").Value = 2.192267813929732
$ws.Range("L12").Value = 0.1507466022469384
$ws.Range("M12").Value = 0.6260057039982172
$ws.Range("O12").Value = 3.756134957206598

$ws.Range("C13").Value = 0.4145498682166817
$ws.Range("D13").Value = 0.190915380071786
$ws.Range("E13").Value = 0.1687717688205872
$ws.Range("F13").Value = 1.543244738710435
$ws.Range("G13").Value = 0.9017696174525582
$ws.Range("H13").Value = 0.9595710528655275
$ws.Range("J13").Value = 0.2018899520821691
$ws.Range("K13").Value = 2.182168304139793
$ws.Range("L13").Value = 0.1507898029978101
$ws.Range("M13").Value = 0.6239785547343502
$ws.Range("O13").Value = 3.75819141408661

$ws.Range("C14").Value = 0.4138384833030955
$ws.Range("D14").Value = 0.1903039053565152
$ws.Range("E14").Value = 0.1687914254786023
$ws.Range("F14").Value = 1.545146761462775
$ws.Range("G14").Value = 0.9033748791404292
$ws.Range("H14").Value = 0.9612046468524369
$ws.Range("J14").Value = 0.2022051017752773
$ws.Range("K14").Value = 2.149223626587798
$ws.Range("L14").Value = 0.1509316636694127
$ws.Range("M14").Value = 0.617367939867691
$ws.Range("O14").Value = 3.764960145020751

$ws.Range("C15").Value = 0.4134062026060406
$ws.Range("D15").Value = 0.1899302795268341
$ws.Range("E15").Value = 0.1688046474786731
$ws.Range("F15").Value = 1.546329772130221
$ws.Range("G15").Value = 0.9043714692432587
$ws.Range("H15").Value = 0.9622146976941366
$ws.Range("J15").Value = 0.2023997548847412
$ws.Range("K15").Value = 2.129037978695067
$ws.Range("L15").Value = 0.1510193058617393
$ws.Range("M15").Value = 0.6133190342795558
$ws.Range("O15").Value = 3.769153957642857

$ws.Range("C16").Value = 0.4109821295916731
$ws.Range("D16").Value = 0.1878033444877047
$ws.Range("E16").Value = 0.1688985504337097
$ws.Range("F16").Value = 1.553384449794038
$ws.Range("G16").Value = 0.9102865297500671
$ws.Range("H16").Value = 0.9681473758140839
$ws.Range("J16").Value = 0.2035399902954289
$ws.Range("K16").Value = 2.013271093080391
$ws.Range("L16").Value = 0.1515330307430514
$ws.Range("M16").Value = 0.5901212319279239
$ws.Range("O16").Value = 3.793918793980822

$ws.Range("C17").Value = 0.4095425479228538
$ws.Range("D17").Value = 0.1865112436413767
$ws.Range("E17").Value = 0.1689723593824155
$ws.Range("F17").Value = 1.557957874496189
$ws.Range("G17").Value = 0.9140970818300431
$ws.Range("H17").Value = 0.9719157967989105
$ws.Range("J17").Value = 0.2042615858920698
$ws.Range("K17").Value = 1.942172878628696
$ws.Range("L17").Value = 0.151858443183837
$ws.Range("M17").Value = 0.5758949303785101
$ws.Range("O17").Value = 3.809763791016991

$ws.Range("C18").Value = 0.408732029318827
$ws.Range("D18").Value = 0.1857727568358314
$ws.Range("E18").Value = 0.1690207754816662
$ws.Range("F18").Value = 1.560678621277482
$ws.Range("G18").Value = 0.9163555682833504
$ws.Range("H18").Value = 0.9741306641507919
$ws.Range("J18").Value = 0.2046847522997943
$ws.Range("K18").Value = 1.901250959708079
$ws.Range("L18").Value = 0.1520493868582768
$ws.Range("M18").Value = 0.5677142928501269
$ws.Range("O18").Value = 3.819117065703892

$ws.Range("C19").Value = 0.4084606086703104
$ws.Range("D19").Value = 0.1855235268451736
$ws.Range("E19").Value = 0.1690381933489817
$ws.Range("F19").Value = 1.561615314623424
$ws.Range("G19").Value = 0.9171317117424493
$ws.Range("H19").Value = 0.9748887158079711
$ws.Range("J19").Value = 0.2048294246885369
$ws.Range("K19").Value = 1.887390761187362
$ws.Range("L19").Value = 0.1521146859435767
$ws.Range("M19").Value = 0.5649448234901797
$ws.Range("O19").Value = 3.822325068351546

$ws.Range("C20").Value = 0.4096939843372525
$ws.Range("D20").Value = 0.1866483046201353
$ws.Range("E20").Value = 0.1689638853288535
$ws.Range("F20").Value = 1.557461686163848
$ws.Range("G20").Value = 0.9136845312347575
$ws.Range("H20").Value = 0.97150973909784
$ws.Range("J20").Value = 0.2041839300974537
$ws.Range("K20").Value = 1.94974433460186
$ws.Range("L20").Value = 0.1518234118826474
$ws.Range("M20").Value = 0.5774091464676445
$ws.Range("O20").Value = 3.808052257236909

$ws.Range("C21").Value = 0.4140467127634793
$ws.Range("D21").Value = 0.1904833187032722
$ws.Range("E21").Value = 0.1687854056188876
$ws.Range("F21").Value = 1.544584346862401
$ws.Range("G21").Value = 0.9029005983893228
$ws.Range("H21").Value = 0.9607228560213485
$ws.Range("J21").Value = 0.2021121981586198
$ws.Range("K21").Value = 2.1589015494406
$ws.Range("L21").Value = 0.1508898398056555
$ws.Range("M21").Value = 0.6193095789638932
$ws.Range("O21").Value = 3.76296204863138

$ws.Range("C22").Value = 0.4170465572915134
$ws.Range("D22").Value = 0.193030974081509
$ws.Range("E22").Value = 0.1687216738040433
$ws.Range("F22").Value = 1.536972342285566
$ws.Range("G22").Value = 0.896448693546219
$ws.Range("H22").Value = 0.9540949381019459
$ws.Range("J22").Value = 0.2008304845330144
$ws.Range("K22").Value = 2.295328146305053
$ws.Range("L22").Value = 0.1503132079495852
$ws.Range("M22").Value = 0.6467072882766871
$ws.Range("O22").Value = 3.735630592048068

$ws.Range("C23").Value = 0.4154313258383695
$ws.Range("D23").Value = 0.1916674949447952
$ws.Range("E23").Value = 0.1687508559403774
$ws.Range("F23").Value = 1.540961412029176
$ws.Range("G23").Value = 0.8998375860792009
$ws.Range("H23").Value = 0.9575937708850546
$ws.Range("J23").Value = 0.2015079480940862
$ws.Range("K23").Value = 2.222540321222539
$ws.Range("L23").Value = 0.150617905980365
$ws.Range("M23").Value = 0.6320835885745311
$ws.Range("O23").Value = 3.750022143567037

$ws.Range("C24").Value = 0.4096254666128516
$ws.Range("D24").Value = 0.1865863257680189
$ws.Range("E24").Value = 0.1689676978121497
$ws.Range("F24").Value = 1.55768572813416
$ws.Range("G24").Value = 0.9138708343403152
$ws.Range("H24").Value = 0.9716931671704288
$ws.Range("J24").Value = 0.2042190123918211
$ws.Range("K24").Value = 1.946321424321809
$ws.Range("L24").Value = 0.1518392375099893
$ws.Range("M24").Value = 0.5767245747716885
$ws.Range("O24").Value = 3.80882528234919

$ws.Range("C25").Value = 0.4040319725911132
$ws.Range("D25").Value = 0.1812852900546318
$ws.Range("E25").Value = 0.1694296875512088
$ws.Range("F25").Value = 1.579175091724984
$ws.Range("G25").Value = 0.9315593200008578
$ws.Range("H25").Value = 0.9887127130486988
$ws.Range("J25").Value = 0.2074539598949272
$ws.Range("K25").Value = 1.647310762050893
$ws.Range("L25").Value = 0.1533010341901679
$ws.Range("M25").Value = 0.5170878799287308
$ws.Range("O25").Value = 3.881414284765441
